# player camera shake, younger and armour consumable
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old heal1/heal2/heal3 rows (rows 2-4) entirely
$ws.Range("A2:C4").ClearContents()

# Rename the "getYounger" consumable to "lowerAge"
$ws.Range("B13").Value = "lowerAge"

# Highlight the armour (row 10) and lowerAge (row 13) consumable rows
$highlight1 = $ws.Range("A10:C10")
$highlight1.Interior.Color = 65535
$highlight1.HorizontalAlignment = -4131

$highlight2 = $ws.Range("A13:C13")
$highlight2.Interior.Color = 65535
$highlight2.HorizontalAlignment = -4131

# Update the selection to the header row range
[void]$ws.Range("A1:F1").Select()
